$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Experimental" property row (A7) gets its Value (B7) set to the literal
# text "false". A plain .Value assignment of "false"/"true" gets coerced to
# a native Excel boolean (t="b"), so instead write it as a text formula and
# collapse it down to a literal value via copy / paste-values -- this keeps
# the cell's shared-string type (t="s") and its original style untouched.
$ws.Range("B7").Formula = "=""false"""
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

# "Case Sensitive" property row (A14) gets its Value (B14) set to "true"
$ws.Range("B14").Formula = "=""true"""
$ws.Range("B14").Copy()
$ws.Range("B14").PasteSpecial(-4163)

$excel.CutCopyMode = 0
